# Update the USD Amount figure on the "Deposit/Crypto/Roobic" row.
# (T1 = "USD Amount" header, T2 holds the day's total that changed
# from 448780 to 449017 in this "Add files via upload" commit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("T2").Value = 449017
